# Add six new "select_one" symptom questions to the survey sheet, inserted
# right after row 168 ("new_section") and before the existing
# who_live_with/who_sharefood_with/who_work_with block, which shifts that
# block (and everything after it) down by six rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# --- 1. Insert six blank rows starting at row 169 -------------------------
$ws.Rows.Item(169).Resize(6).Insert()

# --- 2. Populate the new rows with the new symptom questions --------------
$newRows = @(
    @("select_one yes_no_2", "has_facial_swelling", "Do you have facial swelling?"),
    @("select_one yes_no_2", "has_muscle_fatigue",  "Do you have muscle fatigue?"),
    @("select_one yes_no_2", "has_vomiting",        "Are you vomiting?"),
    @("select_one yes_no_3", "has_cough",            "Do you have a cough?"),
    @("select_one yes_no_4", "has_meningitis",       "Do you have meningitis?"),
    @("select_one yes_no_5", "has_hypertension",     "Do you have hypertension?")
)

$startRow = 169
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $newRows[$i][1]
    $ws.Cells.Item($r, 3).Value = $newRows[$i][2]
}

# --- 3. Widen column B to fit the new, longer question names --------------
$ws.Columns.Item(2).ColumnWidth = 16.8

# --- 4. Update the view: scrolled position & current selection ------------
$ws.Range("C174").Select()
$excel.ActiveWindow.ScrollRow = 161
$excel.ActiveWindow.ScrollColumn = 1
